$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values from 7293 to 7310 for rows 2 through 252
$ws.Range("C2:C252").Value = 7310
